# Add a new "Player Info" sheet before the existing "ODI Batting" sheet,
# and rework the "ODI Batting" sheet's MATCH_CARD_LINK column into a
# MATCH_CODE column holding just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Insert "Player Info" sheet ahead of the existing sheet -------------
$original = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($original)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header styling used elsewhere in the
# workbook.
$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Data row. The ID looks numeric ("5471") but must be stored as text, like
# every other value in this workbook - a leading apostrophe forces text
# entry, then the style is reset so no stray formatting is left behind.
$playerInfo.Range("A2").Value = "'5471"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Ibrahim Zadran"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# --- 2. Update the original "ODI Batting" sheet -----------------------------
# Re-acquire the sheet by name: the handle obtained above now refers to the
# newly inserted sheet since it is positional.
$odi = $wb.Worksheets.Item("ODI Batting")

$odi.Range("D1").Value = "MATCH_CODE"
$odi.Range("D2").Value = "'4379"
$odi.Range("D3").Value = "'4537"
$odi.Range("D4").Value = "'4582"
$odi.Range("D5").Value = "'4585"
$odi.Range("D6").Value = "'4588"
$odi.Range("D7").Value = "'4671"
$odi.Range("D8").Value = "'4674"
$odi.Range("D9").Value = "'4675"
$odi.Range("D2:D9").Style = "Normal"

Write-Output "Player Info sheet added; ODI Batting MATCH_CODE column updated."
